$wb = $excel.ActiveWorkbook

# The shared string "Ready for handoff" is referenced by the Overview sheet's
# B3/C3 cells as well as the per-language sheets' C3 (Status) cells for the
# 4b8ad6e3 handoff row. In the original workbook they all shared the same
# string, so all of them must be updated together to "Handback transform failed".
$newStatus = "Handback transform failed"

# --- Overview sheet: update status for the 4b8ad6e3 handoff row ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# --- zh-cn sheet: update Status + add Error Detail for the 4b8ad6e3 handoff row (row 3) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("K3").Value = "Handback file name: rqoal0o1.5yn is different with handoff file name: 4b8ad6e3-47ad-4821-af0f-38c6eebf06ad.49e6a296e61079771320d1f2174a57e5d4908317.zh-cn."

# --- de-de sheet: update Status + add Error Detail for the 4b8ad6e3 handoff row (row 3) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $newStatus
$dede.Range("K3").Value = "Handback file name: rqoal0o1.5yn is different with handoff file name: 4b8ad6e3-47ad-4821-af0f-38c6eebf06ad.49e6a296e61079771320d1f2174a57e5d4908317.de-de."
